$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated odds/values as described by the commit diff.
# Each line sets a single cell's numeric value to match the target workbook state.
$ws.Range("G3").Value = 3.75
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 2.05
$ws.Range("K3").Value = 2.2
$ws.Range("Y3").Value = 13
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48
$ws.Range("G6").Value = 1.2
$ws.Range("H6").Value = 5.8
$ws.Range("I6").Value = 11.75
$ws.Range("K6").Value = 2.7
$ws.Range("L6").Value = 8.75
$ws.Range("O6").Value = 1.1
$ws.Range("P6").Value = 5.1
$ws.Range("R6").Value = 2.4
$ws.Range("U6").Value = 2
$ws.Range("X6").Value = 6.3
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 7
$ws.Range("AC6").Value = 16.5
$ws.Range("AD6").Value = 12.5
$ws.Range("AH6").Value = 35
$ws.Range("AJ6").Value = 37
$ws.Range("AL6").Value = 150
$ws.Range("AM6").Value = 110
$ws.Range("AN6").Value = 3.1
$ws.Range("AO6").Value = 4.9
$ws.Range("AQ6").Value = 11.25
$ws.Range("AT6").Value = 3.55
$ws.Range("AY6").Value = 65
$ws.Range("BA6").Value = 500
$ws.Range("BB6").Value = 450
$ws.Range("G7").Value = 1.14
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 13
$ws.Range("J7").Value = 1.5
$ws.Range("K7").Value = 3.1
$ws.Range("L7").Value = 10
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 23
$ws.Range("O7").Value = 1.1
$ws.Range("P7").Value = 7
$ws.Range("Q7").Value = 1.36
$ws.Range("R7").Value = 3.1
$ws.Range("S7").Value = 1.2
$ws.Range("T7").Value = 4.33
$ws.Range("U7").Value = 1.91
$ws.Range("V7").Value = 1.8
$ws.Range("W7").Value = 11
$ws.Range("X7").Value = 7.5
$ws.Range("Y7").Value = 11
$ws.Range("Z7").Value = 7.5
$ws.Range("AA7").Value = 11
$ws.Range("AB7").Value = 26
$ws.Range("AC7").Value = 23
$ws.Range("AD7").Value = 15
$ws.Range("AE7").Value = 26
$ws.Range("AF7").Value = 67
$ws.Range("AG7").Value = 251
$ws.Range("AH7").Value = 34
$ws.Range("AI7").Value = 67
$ws.Range("AJ7").Value = 34
$ws.Range("AK7").Value = 151
$ws.Range("AL7").Value = 81
$ws.Range("AM7").Value = 67
$ws.Range("AN7").Value = 3.4
$ws.Range("AO7").Value = 5
$ws.Range("AQ7").Value = 10
$ws.Range("AR7").Value = 29
$ws.Range("AS7").Value = 101
$ws.Range("AT7").Value = 4.33
$ws.Range("AU7").Value = 10
$ws.Range("AX7").Value = 13
$ws.Range("AY7").Value = 51
$ws.Range("AZ7").Value = 41
$ws.Range("BA7").Value = 251
$ws.Range("BB7").Value = 201
$ws.Range("BC7").Value = 301
$ws.Range("G8").Value = 9
$ws.Range("H8").Value = 4.75
$ws.Range("I8").Value = 1.34
$ws.Range("J8").Value = 7.6
$ws.Range("K8").Value = 2.4
$ws.Range("N8").Value = 9
$ws.Range("O8").Value = 1.23
$ws.Range("P8").Value = 3.9
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 2.07
$ws.Range("S8").Value = 1.35
$ws.Range("T8").Value = 3.05
$ws.Range("U8").Value = 2.02
$ws.Range("V8").Value = 1.72
$ws.Range("W8").Value = 18.5
$ws.Range("Y8").Value = 30
$ws.Range("Z8").Value = 300
$ws.Range("AA8").Value = 120
$ws.Range("AB8").Value = 110
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 10
$ws.Range("AG8").Value = 1250
$ws.Range("AH8").Value = 6.4
$ws.Range("AJ8").Value = 9
$ws.Range("AL8").Value = 12
$ws.Range("AM8").Value = 32
$ws.Range("AN8").Value = 9.25
$ws.Range("AO8").Value = 55
$ws.Range("AP8").Value = 50
$ws.Range("AQ8").Value = 450
$ws.Range("AR8").Value = 450
$ws.Range("AT8").Value = 3.05
$ws.Range("AX8").Value = 3.1
$ws.Range("G9").Value = 2.1
$ws.Range("I9").Value = 3.4
$ws.Range("AD9").Value = 6
$ws.Range("AI9").Value = 17
$ws.Range("AJ9").Value = 12
$ws.Range("AK9").Value = 34
$ws.Range("AL9").Value = 26
$ws.Range("AN9").Value = 4.33
$ws.Range("AO9").Value = 12
$ws.Range("BA9").Value = 51
$ws.Range("G12").Value = 9.5
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 12
$ws.Range("Q12").Value = 1.73
$ws.Range("R12").Value = 2.08
$ws.Range("S12").Value = 1.33
$ws.Range("T12").Value = 3.25
$ws.Range("Y12").Value = 23
$ws.Range("AA12").Value = 51
$ws.Range("AC12").Value = 12
$ws.Range("AG12").Value = 401
$ws.Range("AH12").Value = 7
$ws.Range("AI12").Value = 6.5
$ws.Range("AJ12").Value = 9
$ws.Range("AS12").Value = 351
$ws.Range("AT12").Value = 3.25
$ws.Range("AX12").Value = 3.25
$ws.Range("G16").Value = 4.67
$ws.Range("H16").Value = 3.95
$ws.Range("I16").Value = 1.65
$ws.Range("J16").Value = 5.04
$ws.Range("K16").Value = 2.35
$ws.Range("L16").Value = 2.22
$ws.Range("O16").Value = 1.21
$ws.Range("P16").Value = 3.9
$ws.Range("W16").Value = 12
$ws.Range("X16").Value = 23
$ws.Range("Y16").Value = 12
$ws.Range("Z16").Value = 70
$ws.Range("AD16").Value = 6
$ws.Range("AI16").Value = 6.4
$ws.Range("AK16").Value = 9.8
$ws.Range("AM16").Value = 20
$ws.Range("G17").Value = 1.5
$ws.Range("K17").Value = 2.38
$ws.Range("L17").Value = 6
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 13
$ws.Range("O17").Value = 1.22
$ws.Range("P17").Value = 4
$ws.Range("Q17").Value = 1.8
$ws.Range("R17").Value = 2
$ws.Range("AC17").Value = 12
$ws.Range("AE17").Value = 17
$ws.Range("AL17").Value = 41
$ws.Range("AQ17").Value = 21
$ws.Range("AS17").Value = 126
$ws.Range("AU17").Value = 8.5
$ws.Range("BC17").Value = 251
$ws.Range("H18").Value = 7
$ws.Range("L18").Value = 12
$ws.Range("M18").Value = 1.02
$ws.Range("N18").Value = 19
$ws.Range("Q18").Value = 1.48
$ws.Range("R18").Value = 2.6
$ws.Range("S18").Value = 1.25
$ws.Range("T18").Value = 3.75
$ws.Range("X18").Value = 6.5
$ws.Range("AD18").Value = 13
$ws.Range("AH18").Value = 34
$ws.Range("AK18").Value = 251
$ws.Range("AM18").Value = 81
$ws.Range("AT18").Value = 3.75
$ws.Range("AU18").Value = 11
$ws.Range("AX18").Value = 13
$ws.Range("AY18").Value = 51
$ws.Range("BA18").Value = 401
$ws.Range("G19").Value = 3.5
$ws.Range("H19").Value = 3.4
$ws.Range("I19").Value = 2.1
$ws.Range("K19").Value = 2.05
$ws.Range("W19").Value = 9
$ws.Range("AC19").Value = 8.5
$ws.Range("AD19").Value = 6.5
$ws.Range("AG19").Value = 351
$ws.Range("AJ19").Value = 9
$ws.Range("AP19").Value = 29
$ws.Range("AZ19").Value = 23
$ws.Range("Q20").Value = 2.03
$ws.Range("R20").Value = 1.83
$ws.Range("G22").Value = 1.42
$ws.Range("I22").Value = 6.3
$ws.Range("J22").Value = 1.91
$ws.Range("K22").Value = 2.4
$ws.Range("L22").Value = 6
$ws.Range("Q22").Value = 1.65
$ws.Range("R22").Value = 2.12
$ws.Range("T22").Value = 3.1
$ws.Range("W22").Value = 7.4
$ws.Range("Z22").Value = 9.5
$ws.Range("AA22").Value = 11.25
$ws.Range("AE22").Value = 18
$ws.Range("AJ22").Value = 20
$ws.Range("AK22").Value = 120
$ws.Range("AL22").Value = 65
$ws.Range("AM22").Value = 60
$ws.Range("AO22").Value = 6.5
$ws.Range("AP22").Value = 16
$ws.Range("AQ22").Value = 18.5
$ws.Range("AR22").Value = 45
$ws.Range("AT22").Value = 3.1
$ws.Range("AX22").Value = 7.9
$ws.Range("AY22").Value = 37
$ws.Range("AZ22").Value = 37
$ws.Range("BA22").Value = 250
$ws.Range("BB22").Value = 250
$ws.Range("BC22").Value = 450
$ws.Range("Q23").Value = 1.5
$ws.Range("R23").Value = 2.5
$ws.Range("J24").Value = 1.73
$ws.Range("U24").Value = 1.53
$ws.Range("V24").Value = 2.38
$ws.Range("AB24").Value = 19
$ws.Range("AH24").Value = 26
$ws.Range("AU24").Value = 8
$ws.Range("AZ24").Value = 29
$ws.Range("BA24").Value = 101

